$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format cells that hold numeric-looking text so Excel keeps them as text
# (matches the source data which stores price/volume figures as literal strings).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '62.681.66'
$ws.Range("E2").Value = '  +3.38%  '
$ws.Range("D3").Value = '2.445.27'
$ws.Range("E3").Value = '  +2.00%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = '578.67'
$ws.Range("E5").Value = '  +2.99%  '
$ws.Range("D6").Value = '145.55'
$ws.Range("E6").Value = '  +3.14%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("E8").Value = '  +0.56%  '
$ws.Range("D9").Value = '2.443.48'
$ws.Range("E9").Value = '  +1.69%  '
$ws.Range("E10").Value = '  +2.02%  '
$ws.Range("E12").Value = '  +1.52%  '
$ws.Range("E13").Value = '  +3.04%  '
$ws.Range("D14").Value = '28.44'
$ws.Range("E14").Value = '  +9.29%  '
$ws.Range("E15").Value = '  +5.70%  '
$ws.Range("D16").Value = '2.888.90'
$ws.Range("E16").Value = '  +2.12%  '
$ws.Range("D17").Value = '62.561.16'
$ws.Range("E17").Value = '  +3.44%  '
$ws.Range("D18").Value = '2.433.56'
$ws.Range("E18").Value = '  +1.28%  '
$ws.Range("B19").Value = 'Uniswap'
$ws.Range("C19").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D19").Value = '7.77'
$ws.Range("E19").Value = '  -3.52%  '
$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D20").Value = '10.92'
$ws.Range("E20").Value = '  +2.87%  '
$ws.Range("B21").Value = 'BabyDogeCoin'
$ws.Range("C21").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D21").Value = '0.0₆0862'
$ws.Range("E21").Value = '  +209.95%  '
$ws.Range("D22").Value = '326.56'
$ws.Range("E22").Value = '  +1.00%  '
$ws.Range("E24").Value = '  +11.05%  '
$ws.Range("E25").Value = '  +0.05%  '
$ws.Range("D26").Value = '65.52'
$ws.Range("E26").Value = '  +0.98%  '
$ws.Range("D27").Value = '644.37'
$ws.Range("E27").Value = '  +14.69%  '
$ws.Range("D28").Value = '1.17'
$ws.Range("E28").Value = '  +16.65%  '
$ws.Range("D29").Value = '8.47'
$ws.Range("E29").Value = '  +5.71%  '
$ws.Range("D30").Value = '0.0₃0978'
$ws.Range("E30").Value = '  +4.67%  '
$ws.Range("E31").Value = '  +2.04%  '
$ws.Range("D32").Value = '8.17'
$ws.Range("E33").Value = '  +6.43%  '
$ws.Range("E34").Value = '  +3.55%  '
$ws.Range("E35").Value = '  +6.05%  '
$ws.Range("E36").Value = '  +2.26%  '
$ws.Range("D37").Value = '0.999'
$ws.Range("E37").Value = '  +0.11%  '
$ws.Range("E38").Value = '  +3.21%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D39").Value = '5.46'
$ws.Range("E39").Value = '  +6.85%  '
$ws.Range("B40").Value = 'Monero'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D40").Value = '153.20'
$ws.Range("E40").Value = '  +0.00%  '
$ws.Range("E41").Value = '  +0.81%  '
$ws.Range("E42").Value = '  +1.86%  '
$ws.Range("E43").Value = '  +8.83%  '
$ws.Range("D44").Value = '1.76'
$ws.Range("E44").Value = '  +5.47%  '
$ws.Range("D45").Value = '42.57'
$ws.Range("E45").Value = '  +1.96%  '
$ws.Range("E46").Value = '  +0.01%  '
$ws.Range("D47").Value = '15.02'
$ws.Range("E47").Value = '  +28.03%  '
$ws.Range("D48").Value = '144.09'
$ws.Range("E48").Value = '  +1.81%  '
$ws.Range("E49").Value = '  +1.78%  '
$ws.Range("D50").Value = '20.61'
$ws.Range("E50").Value = '  +6.97%  '
$ws.Range("D51").Value = '0.604'
$ws.Range("E51").Value = '  +2.90%  '
